# Atualiza notas dos alunos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new grades in column G (C2) for rows 2-6
$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 1

# Update the active selection to match the author's final cursor position
$ws.Range("G3").Select()
